$wb = $excel.ActiveWorkbook

# --- Remove the old empty "Sheet3" tab -------------------------------------
# Deleting it first frees up sheetId 3, so the sheet we add next reuses it
# (matching the target workbook.xml: the new "testSuite" sheet gets sheetId="3").
$wb.Worksheets.Item("Sheet3").Delete()

# --- Insert a brand-new "testSuite" sheet as the first tab -----------------
$firstSheet = $wb.Worksheets.Item(1)
$testSuite = $wb.Worksheets.Add($firstSheet)
$testSuite.Name = "testSuite"

# Header row first (A1, B1), then the rest of column A, then the rest of
# column B -- this is the write order that reproduces the target
# shared-string table layout.
$testSuite.Range("A1").Value = "TC ID"
$testSuite.Range("B1").Value = "RunMode"

$testSuite.Range("A2").Value = "TestLoginAsBankManager"
$testSuite.Range("A3").Value = "TestAddCustomer"
$testSuite.Range("A4").Value = "TestOpenAccount"

$testSuite.Range("B2").Value = "Y"
$testSuite.Range("B3").Value = "y"
$testSuite.Range("B4").Value = "n"

$testSuite.Columns.Item(1).ColumnWidth = 23.85546875

# --- testOpenAccount: move the selection, make it the non-active tab -------
$openAccount = $wb.Worksheets.Item("testOpenAccount")
$openAccount.Activate()
$openAccount.Range("G11:L21").Select()

# --- testSuite ends up being the active tab/cell, matching the target ------
$testSuite.Activate()
$testSuite.Range("B4").Select()
